# ecg_bom_charger.xlsx edit script
# - fixes incorrect micro USB connector part number (row 13)
# - updates header connector from 1-pos to 2-pos part (row 12)
# - adds Digikey hyperlinks to the Vendor Part Number column
# - updates the saved selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: CONN3 header changes from 1 POS to 2 POS Mill-Max part ---
$ws.Cells.Item(12,3).Value = "ED90503-ND"
$ws.Cells.Item(12,4).Value = "Mill-Max"
$ws.Cells.Item(12,5).Value = "829-22-002-20-001101"
$ws.Cells.Item(12,6).Value = "HEADER 2 POS R/A SIP PCB"
$ws.Cells.Item(12,8).Value = 1
$ws.Cells.Item(12,9).Value = 10.12

# --- Row 13: CONN1,CONN2 corrected micro USB connector part ---
$ws.Cells.Item(13,3).Value = "WM1399CT-ND"
$ws.Cells.Item(13,4).Value = "Molex"
$ws.Cells.Item(13,5).Value = 1050170001
$ws.Cells.Item(13,6).Value = "CONN RCPT MICRO USB R/A SMD"
$ws.Cells.Item(13,9).Value = 1

# --- Add Digikey hyperlinks to the Vendor Part Number column (Column C) ---
# order matters for relationship id assignment - matches original author's order
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.digikey.com/en/products/result?keywords=MCP73832T-2ACI%2FOTCT-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://www.digikey.com/en/products/result?keywords=P470JCT-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://www.digikey.com/en/products/result?keywords=P10KJCT-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://www.digikey.com/en/products/result?keywords=P124KLCT-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.digikey.com/en/products/result?keywords=1276-1056-1-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "https://www.digikey.com/en/products/result?keywords=WM1399CT-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C15"), "https://www.digikey.com/en/products/result?keywords=160-2031-1-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), "https://www.digikey.com/en/products/result?keywords=MF-NSMF150-2CT-ND") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C12"), "https://www.digikey.com/en/products/result?keywords=ED90503-ND") | Out-Null

# --- restore the saved cell selection ---
$ws.Range("D22").Select() | Out-Null
